$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.085.21"
$ws.Range("E2").Value = "  -3.56%  "
$ws.Range("D3").Value = "1.919.18"
$ws.Range("E3").Value = "  -2.73%  "
$ws.Range("E4").Value = "  -0.87%  "
$ws.Range("D5").Value = "'330.81"
$ws.Range("E5").Value = "  +0.55%  "
$ws.Range("E6").Value = "  -0.83%  "
$ws.Range("D7").Value = "'0.4709"
$ws.Range("E7").Value = "  -5.06%  "
$ws.Range("D8").Value = "'0.4038"
$ws.Range("E8").Value = "  -4.06%  "
$ws.Range("D9").Value = "'53.08"
$ws.Range("E9").Value = "  -1.54%  "
$ws.Range("D10").Value = "'0.08421"
$ws.Range("E10").Value = "  -9.34%  "
$ws.Range("E11").Value = "  -5.00%  "
$ws.Range("D12").Value = "'22.17"
$ws.Range("E12").Value = "  -3.08%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.923.50"
$ws.Range("E13").Value = "  -2.55%  "
$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D14").Value = "'7.475"
$ws.Range("E14").Value = "  -5.37%  "
$ws.Range("B15").Value = "Polkadot"
$ws.Range("C15").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D15").Value = "'6.077"
$ws.Range("E15").Value = "  -5.92%  "
$ws.Range("D16").Value = "'1.002"
$ws.Range("E16").Value = "  -1.23%  "
$ws.Range("D17").Value = "'90.07"
$ws.Range("E17").Value = "  -2.05%  "
$ws.Range("D18").Value = "'0.00001063"
$ws.Range("E18").Value = "  -4.36%  "
$ws.Range("D19").Value = "'0.06575"
$ws.Range("E19").Value = "  -2.25%  "
$ws.Range("D20").Value = "'18.09"
$ws.Range("E20").Value = "  -5.82%  "
$ws.Range("E21").Value = "  -0.89%  "
$ws.Range("D22").Value = "'5.729"
$ws.Range("E22").Value = "  -3.94%  "
$ws.Range("D23").Value = "28.059.86"
$ws.Range("E23").Value = "  -3.72%  "
$ws.Range("D24").Value = "'11.35"
$ws.Range("E24").Value = "  -5.30%  "
$ws.Range("D26").Value = "2.088.18"
$ws.Range("E26").Value = "  -5.38%  "
$ws.Range("D27").Value = "'153.93"
$ws.Range("E27").Value = "  -1.28%  "
$ws.Range("D28").Value = "'20.02"
$ws.Range("E28").Value = "  -3.60%  "
$ws.Range("E29").Value = "  -5.55%  "
$ws.Range("E30").Value = "  -8.33%  "
$ws.Range("D31").Value = "'123.45"
$ws.Range("E31").Value = "  -2.99%  "
$ws.Range("D32").Value = "'0.9734"
$ws.Range("E32").Value = "  -7.03%  "
$ws.Range("D33").Value = "'0.09594"
$ws.Range("E33").Value = "  -2.74%  "
$ws.Range("E34").Value = "  -3.60%  "
$ws.Range("D35").Value = "'3.641"
$ws.Range("E35").Value = "  -2.52%  "
$ws.Range("D36").Value = "'5.541"
$ws.Range("E36").Value = "  -4.73%  "
$ws.Range("E37").Value = "  -1.10%  "
$ws.Range("E38").Value = "  -4.86%  "
$ws.Range("D39").Value = "'0.06144"
$ws.Range("E39").Value = "  -4.55%  "
$ws.Range("D40").Value = "'1.218"
$ws.Range("E40").Value = "  -7.99%  "
$ws.Range("D41").Value = "'0.6140"
$ws.Range("E41").Value = "  -5.32%  "
$ws.Range("E42").Value = "  -4.17%  "
$ws.Range("E43").Value = "  -0.82%  "
$ws.Range("E44").Value = "  -5.26%  "
$ws.Range("D45").Value = "'1.308"
$ws.Range("E45").Value = "  -3.63%  "
$ws.Range("D46").Value = "'0.5873"
$ws.Range("E46").Value = "  -5.59%  "
$ws.Range("D47").Value = "'12.74"
$ws.Range("E47").Value = "  -4.27%  "
$ws.Range("D48").Value = "'2.027"
$ws.Range("E48").Value = "  -7.14%  "
$ws.Range("D49").Value = "'3.478"
$ws.Range("E49").Value = "  -0.21%  "
$ws.Range("D50").Value = "'0.06836"
$ws.Range("E50").Value = "  -1.97%  "
$ws.Range("D51").Value = "'110.01"
$ws.Range("E51").Value = "  -2.97%  "
